# Automatische test-sync: 2025-06-19 16:30:10
# Append a new mail-log entry to the "Logs" sheet and refresh the
# "Dashboard" category counts to match.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append row 18 -------------------------------------------
$newRow = 18

$logs.Cells.Item($newRow, 1).Value = "Offerte voor zakelijke samenwerking"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Kunt u mij een offerte sturen voor 100 stuks product X?"
$logs.Cells.Item($newRow, 4).Value = "Bestelling"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 16:28:11"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# --- Dashboard sheet: "Bestelling" now outranks "Klacht" (3 vs 2) --------
$dash.Cells.Item(4, 1).Value = "Bestelling"
$dash.Cells.Item(4, 2).Value = 3
$dash.Cells.Item(5, 1).Value = "Klacht"
$dash.Cells.Item(5, 2).Value = 2

# --- Logs sheet: conditional formatting now spans through row 18 ---------
# Rebuild the two rule groups (Categorie / Beantwoord) against the grown
# ranges, keeping the same rules/order the sheet already had.
$logs.Cells.FormatConditions.Delete()

$catRange = $logs.Range("D2:D18")
$catRange.FormatConditions.Add(1, 3, '="Klacht"') | Out-Null
$catRange.FormatConditions.Add(1, 3, '="Bestelling"') | Out-Null
$catRange.FormatConditions.Add(1, 3, '="Informatieaanvraag"') | Out-Null
$catRange.FormatConditions.Add(1, 3, '="Afmelding"') | Out-Null
$catRange.FormatConditions.Add(1, 3, '="Overig"') | Out-Null
$catRange.FormatConditions.Add(1, 3, '="Retour"') | Out-Null

$answeredRange = $logs.Range("G2:G18")
$answeredRange.FormatConditions.Add(1, 3, '="Ja"') | Out-Null
$answeredRange.FormatConditions.Add(1, 3, '="Nee"') | Out-Null
